$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 44.857143
$ws.Range("I12").Value = 44.857143
$ws.Range("K12").Value = 44.857143
$ws.Range("M12").Value = 125.142857
# Row 53
$ws.Range("H53").Value = 1131.5385
$ws.Range("I53").Value = 1534.1666
$ws.Range("K53").Value = 1534.1666
$ws.Range("M53").Value = -897.1666
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
# Row 132
$ws.Range("H132").Value = 1824.0741
$ws.Range("I132").Value = 1509.7307
$ws.Range("J132").Value = 9997
$ws.Range("K132").Value = 4529.1921
$ws.Range("L132").Value = 29991
$ws.Range("M132").Value = -1999.1921
$ws.Range("N132").Value = -35051
# Row 138
$ws.Range("H138").Value = 2139
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 413.5
$ws.Range("I2").Value = 396.2
$ws.Range("K2").Value = 396.2
$ws.Range("M2").Value = -283.2
# Row 32
$ws.Range("H32").Value = 11888.085
$ws.Range("I32").Value = 11191.536
$ws.Range("J32").Value = 13474.667
$ws.Range("K32").Value = 11191.536
$ws.Range("L32").Value = 13474.667
$ws.Range("M32").Value = -10904.536
$ws.Range("N32").Value = -14048.667
# Row 61
$ws.Range("H61").Value = 874.3333
$ws.Range("I61").Value = 874.3333
$ws.Range("K61").Value = 874.3333
$ws.Range("M61").Value = -662.3333
# Row 74
$ws.Range("H74").Value = 2107.5417
$ws.Range("I74").Value = 1564.8889
$ws.Range("J74").Value = 3735.5
$ws.Range("K74").Value = 1564.8889
$ws.Range("L74").Value = 3735.5
$ws.Range("M74").Value = -690.8888999999999
$ws.Range("N74").Value = -5483.5
# Row 77
$ws.Range("H77").Value = 2107.5417
$ws.Range("I77").Value = 1564.8889
$ws.Range("J77").Value = 3735.5
$ws.Range("K77").Value = 7824.4445
$ws.Range("L77").Value = 18677.5
$ws.Range("M77").Value = -3456.4445
$ws.Range("N77").Value = -27413.5
# Row 116
$ws.Range("H116").Value = 413.5
$ws.Range("I116").Value = 396.2
$ws.Range("K116").Value = 396.2
$ws.Range("M116").Value = 1897.8
# Row 136
$ws.Range("H136").Value = 874.3333
$ws.Range("I136").Value = 874.3333
$ws.Range("K136").Value = 2622.9999
$ws.Range("M136").Value = -72.9998999999998

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 413.5
$ws.Range("I3").Value = 396.2
$ws.Range("K3").Value = 396.2
$ws.Range("M3").Value = -282.2
# Row 86
$ws.Range("H86").Value = 4883.4062
$ws.Range("I86").Value = 3577.5557
$ws.Range("K86").Value = 3577.5557
$ws.Range("M86").Value = -2454.5557
# Row 89
$ws.Range("H89").Value = 4883.4062
$ws.Range("I89").Value = 3577.5557
$ws.Range("K89").Value = 17887.7785
$ws.Range("M89").Value = -12271.7785
# Row 94
$ws.Range("H94").Value = 3341.4092
$ws.Range("I94").Value = 3350.55
$ws.Range("K94").Value = 3350.55
$ws.Range("M94").Value = -2899.55
# Row 107
$ws.Range("H107").Value = 2330.2
$ws.Range("I107").Value = 1821
$ws.Range("J107").Value = 3518.3333
$ws.Range("K107").Value = 1821
$ws.Range("L107").Value = 3518.3333
$ws.Range("M107").Value = 99
$ws.Range("N107").Value = -7358.3333
# Row 113
$ws.Range("H113").Value = 16666666
$ws.Range("I113").Value = 16666666
$ws.Range("K113").Value = 16666666
$ws.Range("M113").Value = -16664496
# Row 134
$ws.Range("H134").Value = 3118.6904
$ws.Range("I134").Value = 2486.4473
$ws.Range("K134").Value = 7459.341899999999
$ws.Range("M134").Value = -4924.341899999999
# Row 140
$ws.Range("H140").Value = 43499.383
$ws.Range("J140").Value = 43499.383
$ws.Range("L140").Value = 43499.383
$ws.Range("N140").Value = -53859.383

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 128
$ws.Range("H128").Value = 378318.84
$ws.Range("I128").Value = 378318.84
$ws.Range("K128").Value = 1134956.52
$ws.Range("M128").Value = -1129976.52
# Row 132
$ws.Range("H132").Value = 6849.7144
$ws.Range("J132").Value = 7583.1665
$ws.Range("L132").Value = 68248.4985
$ws.Range("N132").Value = -73308.4985
# Row 134
$ws.Range("H134").Value = 1984.1818
$ws.Range("I134").Value = 1984.1818
$ws.Range("K134").Value = 5952.5454
$ws.Range("M134").Value = -882.5454
# Row 140
$ws.Range("H140").Value = 2231.6
$ws.Range("I140").Value = 1645.1428
$ws.Range("K140").Value = 4935.428400000001
$ws.Range("M140").Value = 244.5715999999993

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 208422.9
$ws.Range("I122").Value = 323699.28
$ws.Range("J122").Value = 6689.25
$ws.Range("K122").Value = 971097.8400000001
$ws.Range("L122").Value = 20067.75
$ws.Range("M122").Value = -968647.8400000001
$ws.Range("N122").Value = -24967.75

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5054160
$ws.Range("I40").Value = 4037.75
$ws.Range("J40").Value = 18521152
$ws.Range("K40").Value = 4037.75
$ws.Range("L40").Value = 18521152
$ws.Range("M40").Value = -3901.75
$ws.Range("N40").Value = -18521424
# Row 132
$ws.Range("H132").Value = 1936.04
$ws.Range("I132").Value = 1650.45
$ws.Range("K132").Value = 4951.35
$ws.Range("M132").Value = -2421.35
# Row 136
$ws.Range("H136").Value = 5355.3
$ws.Range("I136").Value = 7593
$ws.Range("J136").Value = 1998.75
$ws.Range("K136").Value = 22779
$ws.Range("L136").Value = 5996.25
$ws.Range("M136").Value = -20229
$ws.Range("N136").Value = -11096.25

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 10469.833
$ws.Range("J41").Value = 10519.6
$ws.Range("L41").Value = 10519.6
$ws.Range("N41").Value = -11299.6
# Row 46
$ws.Range("H46").Value = 91520.25
$ws.Range("J46").Value = 91520.25
$ws.Range("L46").Value = 91520.25
$ws.Range("N46").Value = -91982.25
# Row 122
$ws.Range("H122").Value = 2576.6155
$ws.Range("I122").Value = 2262.125
$ws.Range("K122").Value = 6786.375
$ws.Range("M122").Value = -4336.375
# Row 132
$ws.Range("H132").Value = 1225.4348
$ws.Range("I132").Value = 1099.2106
$ws.Range("J132").Value = 1825
$ws.Range("K132").Value = 3297.6318
$ws.Range("L132").Value = 5475
$ws.Range("M132").Value = -767.6318000000001
$ws.Range("N132").Value = -10535
# Row 134
$ws.Range("H134").Value = 91520.25
$ws.Range("J134").Value = 91520.25
$ws.Range("L134").Value = 274560.75
$ws.Range("N134").Value = -279630.75
